$d = $word.ActiveDocument

# --- Activate the built-in Header / Footer paragraph styles (adds style defs to styles.xml) ---
# Use a throw-away paragraph so no explicit pStyle is left behind on real content:
# insert a scratch paragraph at the end, apply the built-in styles to it (which forces
# Word to unhide/emit the style definitions), then delete the scratch paragraph outright.
$scratchRange = $d.Content
$scratchRange.Collapse(0)
$scratchRange.InsertParagraphAfter()
$scratchPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$scratchPara.Range.Style = $d.Styles.Item(-32)
$scratchPara.Range.Style = $d.Styles.Item(-33)
$scratchPara.Range.Delete()

# --- Text replacements: "Info center(re)" -> "centre d'informations" ---
# Paragraph 1: "Info centre" -> "centre d'informations"
$d.Content.Find.Execute("Info centre", $true, $false, $false, $false, $false, $true, 1, $false, "centre d’informations", 2)

# Paragraph 2, occurrence 1: "l'info center" -> "le centre d'informations"
$d.Content.Find.Execute("l’info center", $true, $false, $false, $false, $false, $true, 1, $false, "le centre d’informations", 2)

# Paragraph 2, occurrence 2: "accéder l'Info center" -> "accéder au centre d'informations"
$d.Content.Find.Execute("accéder l’Info center", $true, $false, $false, $false, $false, $true, 1, $false, "accéder au centre d’informations", 2)

# Paragraph 2, occurrence 3: "grâce à l'Info center" -> "grâce au centre d'informations"
$d.Content.Find.Execute("grâce à l’Info center", $true, $false, $false, $false, $false, $true, 1, $false, "grâce au centre d’informations", 2)

# --- Paragraph justification: set "both" (justify) alignment on every paragraph ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $d.Paragraphs.Item($i).Range.ParagraphFormat.Alignment = 3
}

# --- Add a second empty justified paragraph after the current trailing empty paragraph ---
# (Assigning a bare paragraph-mark character to a collapsed Range at the very end of the
# story -- rather than InsertParagraphAfter -- avoids leaving a stray empty <w:r/> behind.)
$endR = $d.Range($d.Content.End, $d.Content.End)
$endR.Text = [char]13

# Re-apply justification to all paragraphs (including the newly inserted one)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $d.Paragraphs.Item($i).Range.ParagraphFormat.Alignment = 3
}

# --- Materialize footnotes.xml / endnotes.xml (separator + continuation separator) ---
$fnRange = $d.Paragraphs.Item(1).Range
$fnRange.Collapse(0)
$fn = $d.Footnotes.Add($fnRange, "", "x")
$fn.Delete()

$enRange = $d.Paragraphs.Item(1).Range
$enRange.Collapse(0)
$en = $d.Endnotes.Add($enRange, "", "x")
$en.Delete()
